$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D values remain text (avoid Excel auto-converting numeric-looking strings)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '54.141.88'
$ws.Range('E2').Value = '  +0.20%  '
$ws.Range('D3').Value = '2.289.87'
$ws.Range('E3').Value = '  +1.35%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '495.73'
$ws.Range('E5').Value = '  +1.28%  '
$ws.Range('D6').Value = '127.62'
$ws.Range('E6').Value = '  +0.70%  '
$ws.Range('D7').Value = '0.996'
$ws.Range('E7').Value = '  -0.30%  '
$ws.Range('E8').Value = '  +1.80%  '
$ws.Range('D9').Value = '2.287.65'
$ws.Range('E9').Value = '  +1.00%  '
$ws.Range('D10').Value = '0.0949'
$ws.Range('E10').Value = '  +3.04%  '
$ws.Range('E11').Value = '  +2.46%  '
$ws.Range('E12').Value = '  +3.10%  '
$ws.Range('E13').Value = '  -2.42%  '
$ws.Range('D14').Value = '2.692.76'
$ws.Range('E14').Value = '  +1.21%  '
$ws.Range('D15').Value = '21.78'
$ws.Range('E15').Value = '  +3.57%  '
$ws.Range('D16').Value = '54.267.56'
$ws.Range('E16').Value = '  +0.53%  '
$ws.Range('E17').Value = '  +0.60%  '
$ws.Range('D18').Value = '2.270.05'
$ws.Range('E18').Value = '  +1.08%  '
$ws.Range('D19').Value = '10.08'
$ws.Range('E19').Value = '  +4.84%  '
$ws.Range('D20').Value = '4.12'
$ws.Range('E20').Value = '  +3.65%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').Value = '302.12'
$ws.Range('E21').Value = '  +0.48%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').Value = '6.47'
$ws.Range('E22').Value = '  +5.76%  '
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  -0.22%  '
$ws.Range('E24').Value = '  -1.57%  '
$ws.Range('D25').Value = '62.48'
$ws.Range('E25').Value = '  -1.87%  '
$ws.Range('D26').Value = '0.998'
$ws.Range('E26').Value = '  -0.15%  '
$ws.Range('E27').Value = '  +2.05%  '
$ws.Range('E28').Value = '  +5.40%  '
$ws.Range('D29').Value = '2.386.39'
$ws.Range('E29').Value = '  +0.77%  '
$ws.Range('D30').Value = '7.08'
$ws.Range('E30').Value = '  +0.16%  '
$ws.Range('D31').Value = '169.20'
$ws.Range('E31').Value = '  -0.19%  '
$ws.Range('B32').Value = 'PEPE'
$ws.Range('C32').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D32').Value = '0.0₃0688'
$ws.Range('E32').Value = '  -0.19%  '
$ws.Range('B33').Value = 'PancakeSwap'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').Value = '1.60'
$ws.Range('E33').Value = '  +0.07%  '
$ws.Range('D34').Value = '5.89'
$ws.Range('E34').Value = '  +2.62%  '
$ws.Range('D35').Value = '0.997'
$ws.Range('E35').Value = '  -0.09%  '
$ws.Range('D36').Value = '0.996'
$ws.Range('E36').Value = '  -0.19%  '
$ws.Range('E37').Value = '  +0.78%  '
$ws.Range('D38').Value = '17.63'
$ws.Range('E38').Value = '  +1.28%  '
$ws.Range('E39').Value = '  +2.82%  '
$ws.Range('E40').Value = '  +4.09%  '
$ws.Range('D41').Value = '3.74'
$ws.Range('E41').Value = '  +4.19%  '
$ws.Range('D42').Value = '35.41'
$ws.Range('E43').Value = '  +2.70%  '
$ws.Range('E44').Value = '  +1.94%  '
$ws.Range('E45').Value = '  +1.59%  '
$ws.Range('B46').Value = 'RenderToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D46').Value = '4.96'
$ws.Range('E46').Value = '  +6.28%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').Value = '127.76'
$ws.Range('E47').Value = '  +4.49%  '
$ws.Range('D48').Value = '0.0890'
$ws.Range('E48').Value = '  +1.48%  '
$ws.Range('E49').Value = '  +1.13%  '
$ws.Range('D50').Value = '240.04'
$ws.Range('E50').Value = '  +1.00%  '
$ws.Range('D51').Value = '0.0485'
$ws.Range('E51').Value = '  +2.76%  '
